$d = $word.ActiveDocument
$sec = $d.Sections(1)

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2, wdHeaderFooterEvenPages = 3

# "First page" header (header1.xml) holds the BTec logo: image2.jpg -> image1.jpg
$hdrFirst = $sec.Headers(2)
if ($hdrFirst.Range.InlineShapes.Count -gt 0) {
    $btecLogo = $hdrFirst.Range.InlineShapes(1)
    $btecLogo.Name = "image1.jpg"
}

# "First page" footer (footer1.xml) holds a Pearson logo: image1.png -> image2.png
$ftrFirst = $sec.Footers(1)
if ($ftrFirst.Range.InlineShapes.Count -gt 0) {
    $pearsonLogo1 = $ftrFirst.Range.InlineShapes(1)
    $pearsonLogo1.Name = "image2.png"
}

# "Primary/default" footer (footer2.xml) holds a Pearson logo: image1.png -> image2.png
$ftrPrimary = $sec.Footers(2)
if ($ftrPrimary.Range.InlineShapes.Count -gt 0) {
    $pearsonLogo2 = $ftrPrimary.Range.InlineShapes(1)
    $pearsonLogo2.Name = "image2.png"
}
